# Adds the RUC/CI, DIRECCION, TELEFONO y CORREO columns to the client
# upload template (module de configuracion / generacion de pdf_guia),
# expanding the header band from column A to column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clone the existing header-band formatting onto the new columns ---
# Row 1 (title fill/merge style) and row 2 (border strip) keep column A's
# look; row 3 (column headers) keeps the same style used by the existing
# "NOMBRE" header cell.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("B1:E1").PasteSpecial(-4122) | Out-Null

$ws.Range("A2").Copy() | Out-Null
$ws.Range("B2:E2").PasteSpecial(-4122) | Out-Null

$ws.Range("A3").Copy() | Out-Null
$ws.Range("B3:E3").PasteSpecial(-4122) | Out-Null

# --- Re-merge the "CLIENTES" title band across the new columns ---
$ws.Range("A1:E2").Merge() | Out-Null

# --- New header row: RUC/CI | NOMBRE | DIRECCION | TELEFONO | CORREO ---
$ws.Range("A3").Value = "RUC/CI"
$ws.Range("B3").Value = "NOMBRE"
$ws.Range("C3").Value = "DIRECCION"
$ws.Range("D3").Value = "TELEFONO"
$ws.Range("E3").Value = "CORREO"

# --- Column widths sized to fit the new header captions ---
$ws.Columns("A").ColumnWidth = 12.833333333333334
$ws.Columns("B").ColumnWidth = 24.666666666666668
$ws.Columns("C").ColumnWidth = 19.833333333333332
$ws.Columns("D").ColumnWidth = 21.666666666666668
$ws.Columns("E").ColumnWidth = 22.666666666666668

# --- Reset the active selection to the first data row ---
$ws.Range("A4").Select() | Out-Null
